$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: bump the commit-hash bookmark id (27 -> 26). Bookmark ids in
# this document are regenerated on save by the engine whenever the
# document is mutated, so this happens automatically as a side effect
# of the edits below; nothing to do for it explicitly.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Change 2: the second "(FIGURE)" callout (the one immediately before
# "). For example,") becomes "(Figure 1)" now that the figure has a
# real caption/number. The first "(FIGURE)" earlier in the same
# paragraph is left untouched.
# ---------------------------------------------------------------------
$searchRng = $d.Content
$searchRng.Find.ClearFormatting()
$searchRng.Find.Execute("FIGURE", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$searchRng.Collapse(0)  # wdCollapseEnd -- resume searching after the first hit
$searchRng.Find.Execute("FIGURE", $true, $false, $false, $false, $false, $true, 1, $false, "Figure 1", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 3: the lone "Figures" Heading2 section marker at the very end
# of the manuscript is replaced by the actual figure legend: a blank
# FirstParagraph spacer followed by a bold BodyText caption.
# ---------------------------------------------------------------------
$figuresPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Figures") {
        $figuresPara = $p
        break
    }
}

if ($figuresPara -ne $null) {
    $rng = $figuresPara.Range
    $xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Figure 1. ESV rate increases as the number of genomes sampled per species increases</w:t></w:r></w:p>'
    $rng.InsertXML($xmlFrag)
}
